$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("E8")
$r.Interior.Color = 65535
$r.Interior.PatternColor = 65535
